$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new ensemble prediction values for row 35 (D35:I35)
$ws.Range("D35").Value = 0.9004
$ws.Range("E35").Value = 0.6717
$ws.Range("F35").Value = 5
$ws.Range("G35").Value = 0.9009
$ws.Range("H35").Value = 0.7082
$ws.Range("I35").Value = 6

# Match the formatting of the already-populated C35 cell (center aligned, General format)
$ws.Range("D35:I35").HorizontalAlignment = -4108

# Update the view: scroll position and active selection
$ws.Range("I35").Select()
$excel.ActiveWindow.ScrollRow = 13
